$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking strings stored as text (inlineStr) in the
# source workbook. Force the cell to Text format before assigning so the COM layer
# keeps the value as a string instead of silently parsing it into a Double, then
# restore the default "Normal" style so no visible formatting change is introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '243.68'
Set-TextValue 'D3' '23.77'
Set-TextValue 'D4' '5.257'
Set-TextValue 'D5' '0.05814'
Set-TextValue 'D6' '6.474'
Set-TextValue 'D7' '3.335'
Set-TextValue 'D8' '0.8079'
Set-TextValue 'D9' '0.8735'
Set-TextValue 'D10' '0.1382'
Set-TextValue 'D11' '0.07257'
Set-TextValue 'D12' '0.03076'
Set-TextValue 'D13' '0.03052'
Set-TextValue 'D14' '0.09320'
Set-TextValue 'D15' '3.863'
Set-TextValue 'D16' '0.001542'
Set-TextValue 'D17' '0.04692'
Set-TextValue 'D18' '0.0006048'
Set-TextValue 'D19' '0.006157'
Set-TextValue 'D21' '0.004594'
Set-TextValue 'D22' '0.00008697'
Set-TextValue 'D24' '2.175'
Set-TextValue 'D25' '0.3200'
Set-TextValue 'D26' '0.1318'
Set-TextValue 'D40' '0.03779'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006318'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1053'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.002614'
$ws.Range("E43").Value = '42CEJICEJI'
Set-TextValue 'D44' '0.007977'
Set-TextValue 'D45' '0.00005522'
Set-TextValue 'D47' '0.5498'
Set-TextValue 'D48' '0.01413'
